$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("wumpus", 0.01, 0.01, 58.09999999999992, "0.182"),
    @("wumpus", 0.01, 0.1,  58.09999999999992, "0.182"),
    @("wumpus", 0.1,  0.9,  58.09999999999992, "0.182"),
    @("wumpus", 0.1,  0.5,  58.09999999999992, "0.182"),
    @("wumpus", 0.1,  0.01, 65.20874699999997, "0.204"),
    @("wumpus", 0.01, 0.5,  0,                 "0.000"),
    @("wumpus", 0.01, 0.9,  0,                 "0.000"),
    @("wumpus", 0.1,  0.1,  64.52430000000004, "0.202"),
    @("wumpus", 0.5,  0.1,  161.46824210875,   "0.505"),
    @("wumpus", 0.5,  0.5,  151.7745720625,    "0.474"),
    @("wumpus", 0.5,  0.01, 165.5361683904909, "0.517"),
    @("wumpus", 0.5,  0.9,  141.6811737500001, "0.443"),
    @("wumpus", 0.9,  0.9,  3011.463534858126, "9.411"),
    @("wumpus", 0.9,  0.01, 3050.369902380611, "9.532"),
    @("wumpus", 0.9,  0.5,  3028.382329004643, "9.464"),
    @("wumpus", 0.9,  0.1,  3046.927929267263, "9.522")
)

$startRow = 18
$endRow = $startRow + $data.Count - 1

# Ensure column E is treated as text so values like "0.182" are not
# auto-converted to numbers.
$ws.Range("E$startRow" + ":E$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
